# Scheduled price-refresh: updates currentAveragePrice / LevePrice / LeveProfit
# columns (H:N) with freshly pulled market-board data across the Leve worksheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2259.64
$ws.Range("I40").Value = 1483.3334
$ws.Range("J40").Value = 2365.5
$ws.Range("K40").Value = 1483.3334
$ws.Range("L40").Value = 2365.5
$ws.Range("M40").Value = -1308.3334
$ws.Range("N40").Value = -2715.5
$ws.Range("H92").Value = 1985.1333
$ws.Range("I92").Value = 2210.5833
$ws.Range("J92").Value = 1083.3334
$ws.Range("K92").Value = 2210.5833
$ws.Range("L92").Value = 1083.3334
$ws.Range("M92").Value = -962.5832999999998
$ws.Range("N92").Value = -3579.3334
$ws.Range("H99").Value = 191.25
$ws.Range("I99").Value = 191.25
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 573.75
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 924.25
$ws.Range("N99").ClearContents()
$ws.Range("H100").Value = 1660.7
$ws.Range("I100").Value = 1534.1111
$ws.Range("J100").Value = 2800
$ws.Range("K100").Value = 1534.1111
$ws.Range("L100").Value = 2800
$ws.Range("M100").Value = -993.1111000000001
$ws.Range("N100").Value = -3882
$ws.Range("H101").Value = 111670.22
$ws.Range("I101").Value = 436
$ws.Range("J101").Value = 500990
$ws.Range("K101").Value = 1308
$ws.Range("L101").Value = 1502970
$ws.Range("M101").Value = 314
$ws.Range("N101").Value = -1506214
$ws.Range("H106").Value = 3586.2856
$ws.Range("I106").Value = 3620.8
$ws.Range("K106").Value = 3620.8
$ws.Range("M106").Value = -2989.8
$ws.Range("H139").Value = 85320.25
$ws.Range("J139").Value = 85320.25
$ws.Range("L139").Value = 85320.25
$ws.Range("N139").Value = -95600.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5365.609
$ws.Range("I32").Value = 4402.6
$ws.Range("K32").Value = 4402.6
$ws.Range("M32").Value = -4115.6
$ws.Range("H61").Value = 4000
$ws.Range("I61").Value = 5000
$ws.Range("K61").Value = 5000
$ws.Range("M61").Value = -4788
$ws.Range("H62").Value = 31900
$ws.Range("J62").Value = 31900
$ws.Range("L62").Value = 31900
$ws.Range("N62").Value = -33148
$ws.Range("H63").Value = 1990
$ws.Range("I63").Value = 1990
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 1990
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1304
$ws.Range("N63").ClearContents()
$ws.Range("H64").Value = 28333.334
$ws.Range("J64").Value = 28333.334
$ws.Range("L64").Value = 28333.334
$ws.Range("N64").Value = -28829.334
$ws.Range("H65").Value = 31900
$ws.Range("J65").Value = 31900
$ws.Range("L65").Value = 95700
$ws.Range("N65").Value = -101940
$ws.Range("H66").Value = 1990
$ws.Range("I66").Value = 1990
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 9950
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -6518
$ws.Range("N66").ClearContents()
$ws.Range("H67").Value = 28333.334
$ws.Range("J67").Value = 28333.334
$ws.Range("L67").Value = 28333.334
$ws.Range("N67").Value = -30049.334
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H74").Value = 1018.9167
$ws.Range("I74").Value = 950.26666
$ws.Range("J74").Value = 1133.3334
$ws.Range("K74").Value = 950.26666
$ws.Range("L74").Value = 1133.3334
$ws.Range("M74").Value = -76.26666
$ws.Range("N74").Value = -2881.3334
$ws.Range("H77").Value = 1018.9167
$ws.Range("I77").Value = 950.26666
$ws.Range("J77").Value = 1133.3334
$ws.Range("K77").Value = 4751.3333
$ws.Range("L77").Value = 5666.666999999999
$ws.Range("M77").Value = -383.3333000000002
$ws.Range("N77").Value = -14402.667
$ws.Range("H97").Value = 676.7143
$ws.Range("I97").Value = 796
$ws.Range("J97").Value = 462
$ws.Range("K97").Value = 796
$ws.Range("L97").Value = 462
$ws.Range("M97").Value = -300
$ws.Range("N97").Value = -1454
$ws.Range("H102").Value = 4385.067
$ws.Range("I102").Value = 4342.273
$ws.Range("K102").Value = 4342.273
$ws.Range("M102").Value = -2720.273
$ws.Range("H132").Value = 2158.0312
$ws.Range("I132").Value = 1761.4814
$ws.Range("J132").Value = 4299.4
$ws.Range("K132").Value = 5284.4442
$ws.Range("L132").Value = 12898.2
$ws.Range("M132").Value = -2754.4442
$ws.Range("N132").Value = -17958.2
$ws.Range("H136").Value = 4000
$ws.Range("I136").Value = 5000
$ws.Range("K136").Value = 15000
$ws.Range("M136").Value = -12450
$ws.Range("H139").Value = 35000
$ws.Range("J139").Value = 35000
$ws.Range("L139").Value = 35000
$ws.Range("N139").Value = -45280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 906.4375
$ws.Range("I94").Value = 626
$ws.Range("J94").Value = 1267
$ws.Range("K94").Value = 626
$ws.Range("L94").Value = 1267
$ws.Range("M94").Value = -175
$ws.Range("N94").Value = -2169
$ws.Range("H99").Value = 25642592
$ws.Range("I99").Value = 35715572
$ws.Range("J99").Value = 2276.7273
$ws.Range("K99").Value = 35715572
$ws.Range("L99").Value = 2276.7273
$ws.Range("M99").Value = -35714074
$ws.Range("N99").Value = -5272.7273
$ws.Range("H107").Value = 9321.375
$ws.Range("I107").Value = 1095.8572
$ws.Range("J107").Value = 66900
$ws.Range("K107").Value = 1095.8572
$ws.Range("L107").Value = 66900
$ws.Range("M107").Value = 824.1428000000001
$ws.Range("N107").Value = -70740
$ws.Range("H134").Value = 94009.09
$ws.Range("I134").Value = 5785.7144
$ws.Range("K134").Value = 17357.1432
$ws.Range("M134").Value = -14822.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 5000
$ws.Range("J99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("N99").ClearContents()
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").ClearContents()
$ws.Range("H115").Value = 1539.6
$ws.Range("I115").Value = 954
$ws.Range("K115").Value = 2862
$ws.Range("M115").Value = -1687
$ws.Range("H121").Value = 1754.2858
$ws.Range("I121").Value = 315
$ws.Range("J121").Value = 2330
$ws.Range("K121").Value = 945
$ws.Range("L121").Value = 6990
$ws.Range("M121").Value = 365
$ws.Range("N121").Value = -9610

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1520.4445
$ws.Range("I113").Value = 1398
$ws.Range("J113").Value = 1581.6666
$ws.Range("K113").Value = 1398
$ws.Range("L113").Value = 1581.6666
$ws.Range("M113").Value = 772
$ws.Range("N113").Value = -5921.6666
$ws.Range("H126").Value = 2659.8572
$ws.Range("I126").Value = 2913.7778
$ws.Range("J126").Value = 2202.8
$ws.Range("K126").Value = 8741.3334
$ws.Range("L126").Value = 6608.400000000001
$ws.Range("M126").Value = -6271.3334
$ws.Range("N126").Value = -11548.4
$ws.Range("H132").Value = 3681.5557
$ws.Range("I132").Value = 3037
$ws.Range("J132").Value = 4970.6665
$ws.Range("K132").Value = 9111
$ws.Range("L132").Value = 14911.9995
$ws.Range("M132").Value = -6581
$ws.Range("N132").Value = -19971.9995
$ws.Range("H138").Value = 34114.5
$ws.Range("J138").Value = 34114.5
$ws.Range("L138").Value = 34114.5
$ws.Range("N138").Value = -44394.5
$ws.Range("H139").Value = 23295.334
$ws.Range("J139").Value = 23295.334
$ws.Range("L139").Value = 23295.334
$ws.Range("N139").Value = -33575.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 12000
$ws.Range("J6").Value = 12000
$ws.Range("L6").Value = 12000
$ws.Range("N6").Value = -12224
$ws.Range("H16").Value = 3533.3333
$ws.Range("I16").Value = 10000
$ws.Range("K16").Value = 10000
$ws.Range("M16").Value = -9830
$ws.Range("H93").Value = 1255.5927
$ws.Range("I93").Value = 945.05
$ws.Range("K93").Value = 945.05
$ws.Range("M93").Value = 302.95
$ws.Range("H138").Value = 29132.25
$ws.Range("J138").Value = 29132.25
$ws.Range("L138").Value = 29132.25
$ws.Range("N138").Value = -39412.25
$ws.Range("H139").Value = 57310
$ws.Range("J139").Value = 57310
$ws.Range("L139").Value = 57310
$ws.Range("N139").Value = -67590

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 523.6
$ws.Range("I113").Value = 490.75
$ws.Range("J113").Value = 655
$ws.Range("K113").Value = 1472.25
$ws.Range("L113").Value = 1965
$ws.Range("M113").Value = 697.75
$ws.Range("N113").Value = -6305
$ws.Range("H122").Value = 2460.2856
$ws.Range("I122").Value = 2456.5
$ws.Range("J122").Value = 2467.8572
$ws.Range("K122").Value = 7369.5
$ws.Range("L122").Value = 7403.571599999999
$ws.Range("M122").Value = -4919.5
$ws.Range("N122").Value = -12303.5716
$ws.Range("H126").Value = 617.4
$ws.Range("I126").Value = 608.2222
$ws.Range("J126").Value = 700
$ws.Range("K126").Value = 1824.6666
$ws.Range("L126").Value = 2100
$ws.Range("M126").Value = 645.3334
$ws.Range("N126").Value = -7040
$ws.Range("H138").Value = 42143
$ws.Range("J138").Value = 42143
$ws.Range("L138").Value = 42143
$ws.Range("N138").Value = -52423
